$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (date style) of the last existing data row (185)
# down into the two new rows so the new date cells pick up the same
# "yyyy-mm-dd hh:mm:ss" style used throughout column A.
$ws.Range("A185").Copy()
$ws.Range("A186:A187").PasteSpecial(-4122)

# ---- Row 186 ----
$ws.Range("A186").Value = 45502.2916666667
$ws.Range("B186").Value = 0
$ws.Range("C186").Value = 2.75999999046326
$ws.Range("D186").Value = 2.75999999046326
$ws.Range("E186").Value = 2.75999999046326
$ws.Range("F186").Value = 2.75999999046326
$ws.Range("H186").Value = "XHS.MI"

# ---- Row 187 ----
$ws.Range("A187").Value = 45503.5047453704
$ws.Range("B187").Value = 1500
$ws.Range("C187").Value = 2.85999989509583
$ws.Range("D187").Value = 2.83999991416931
$ws.Range("E187").Value = 2.83999991416931
$ws.Range("F187").Value = 2.85999989509583
$ws.Range("H187").Value = "XHS.MI"

# ---- Column G (adj_close) must be stored as TEXT (matches source data,
# which keeps these numeric-looking strings as shared strings, not
# numbers). Using a scratch formula + paste-values round-trip forces the
# text type without leaving behind a stray number-format style. ----
$ws.Range("Z1").Formula = "=""2.75999999046326"""
$ws.Range("Z1").Copy()
$ws.Range("G186").PasteSpecial(-4163)
$ws.Range("Z1").ClearContents()

$ws.Range("Z1").Formula = "=""2.85999989509583"""
$ws.Range("Z1").Copy()
$ws.Range("G187").PasteSpecial(-4163)
$ws.Range("Z1").ClearContents()
